$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.944.50"
$ws.Range("E2").Value = "  +0.17%  "
$ws.Range("D3").Value = "2.796.74"
$ws.Range("E3").Value = "  -0.77%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'359.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.31%  "
$ws.Range("D6").Value = "'109.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  -0.58%  "
$ws.Range("E9").Value = "  -1.14%  "
$ws.Range("D10").Value = "'40.27"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.19%  "
$ws.Range("D11").Value = "'0.0853"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.16%  "
$ws.Range("E12").Value = "  +2.20%  "
$ws.Range("D13").Value = "'19.58"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.62%  "
$ws.Range("D14").Value = "'7.65"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.28%  "
$ws.Range("D15").Value = "3.231.04"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "2.798.69"
$ws.Range("E16").Value = "  -0.48%  "
$ws.Range("D17").Value = "'0.944"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.84%  "
$ws.Range("D18").Value = "51.901.35"
$ws.Range("E18").Value = "  +0.23%  "
$ws.Range("D19").Value = "'7.46"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.01%  "
$ws.Range("E20").Value = "  -1.94%  "
$ws.Range("D21").Value = "'13.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("E22").Value = "  -1.03%  "
$ws.Range("D23").Value = "'70.29"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.58%  "
$ws.Range("D24").Value = "'270.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.21%  "
$ws.Range("E25").Value = "  -0.64%  "
$ws.Range("D26").Value = "'26.53"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.64%  "
$ws.Range("E27").Value = "  +0.11%  "
$ws.Range("D28").Value = "'0.165"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +18.52%  "
$ws.Range("D29").Value = "'10.29"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.02%  "
$ws.Range("D30").Value = "'2.17"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.87%  "
$ws.Range("D31").Value = "'6.16"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.74%  "
$ws.Range("E32").Value = "  -1.38%  "
$ws.Range("D33").Value = "'34.97"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.78%  "
$ws.Range("D34").Value = "'0.0467"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.05%  "
$ws.Range("E35").Value = "  +0.69%  "
$ws.Range("E36").Value = "  -3.30%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("E38").Value = "  +2.35%  "
$ws.Range("D39").Value = "'3.22"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.78%  "
$ws.Range("E40").Value = "  -3.10%  "
$ws.Range("D41").Value = "'2.61"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.34%  "
$ws.Range("E42").Value = "  -1.57%  "
$ws.Range("E43").Value = "  -1.46%  "
$ws.Range("D44").Value = "'119.24"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -3.95%  "
$ws.Range("D45").Value = "'21.91"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.23%  "
$ws.Range("D46").Value = "2.090.20"
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("E47").Value = "  -1.48%  "
$ws.Range("D48").Value = "'2.27"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.59%  "
$ws.Range("E49").Value = "  -3.55%  "
$ws.Range("D50").Value = "'0.948"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.98%  "
$ws.Range("E51").Value = "  +30.45%  "
